# Update Sheets via scheduled runner
# Applies updated currentAveragePrice / LevePrice / LeveProfit figures
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 300
$ws.Range("I5").Value = 300
$ws.Range("K5").Value = 300
$ws.Range("M5").Value = -185
$ws.Range("H18").Value = 12737.625
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H70").Value = 9678.277
$ws.Range("J70").Value = 10756.8125
$ws.Range("L70").Value = 32270.4375
$ws.Range("N70").Value = -32810.4375
$ws.Range("H73").Value = 9678.277
$ws.Range("J73").Value = 10756.8125
$ws.Range("L73").Value = 32270.4375
$ws.Range("N73").Value = -34142.4375
$ws.Range("H132").Value = 3635.6667
$ws.Range("I132").Value = 3155.0476
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 9465.1428
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -6935.1428
$ws.Range("N132").Value = -26060

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 40000
$ws.Range("J119").Value = 40000
$ws.Range("L119").Value = 40000
$ws.Range("N119").Value = -49676
$ws.Range("H123").Value = 48000
$ws.Range("J123").Value = 48000
$ws.Range("L123").Value = 48000
$ws.Range("N123").Value = -57800
$ws.Range("H130").Value = 100000
$ws.Range("J130").Value = 100000
$ws.Range("L130").Value = 100000
$ws.Range("N130").Value = -110040
$ws.Range("H131").Value = 100000
$ws.Range("J131").Value = 100000
$ws.Range("L131").Value = 100000
$ws.Range("N131").Value = -110080

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H20").Value = 893.3333
$ws.Range("I20").Value = 340
$ws.Range("K20").Value = 340
$ws.Range("M20").Value = -93
$ws.Range("H105").Value = 4604
$ws.Range("I105").Value = 4340
$ws.Range("K105").Value = 4340
$ws.Range("M105").Value = -2593

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 419.8
$ws.Range("J12").Value = 262
$ws.Range("L12").Value = 262
$ws.Range("N12").Value = -602
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H22").Value = 1000.1667
$ws.Range("I22").Value = 750.25
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 750.25
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -400.25
$ws.Range("N22").Value = -2200
$ws.Range("H25").Value = 9000
$ws.Range("I25").Value = 9000
$ws.Range("K25").Value = 9000
$ws.Range("M25").Value = -8826
$ws.Range("H94").Value = 221.5
$ws.Range("J94").Value = 221.5
$ws.Range("L94").Value = 221.5
$ws.Range("N94").Value = -1123.5
$ws.Range("H98").Value = 200000
$ws.Range("I98").Value = 200000
$ws.Range("K98").Value = 200000
$ws.Range("M98").Value = -197754
$ws.Range("H130").Value = 49875
$ws.Range("J130").Value = 49875
$ws.Range("L130").Value = 49875
$ws.Range("N130").Value = -59915
$ws.Range("H132").Value = 1966.6666
$ws.Range("I132").Value = 1450
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 4350
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -1820
$ws.Range("N132").Value = -14060

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 500
$ws.Range("I19").Value = 500
$ws.Range("K19").Value = 1500
$ws.Range("M19").Value = -1326
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H35").Value = 600
$ws.Range("J35").Value = 600
$ws.Range("L35").Value = 1800
$ws.Range("N35").Value = -2376
$ws.Range("H44").Value = 311.2
$ws.Range("I44").Value = 102
$ws.Range("J44").Value = 625
$ws.Range("K44").Value = 306
$ws.Range("L44").Value = 1875
$ws.Range("M44").Value = 92
$ws.Range("N44").Value = -2671
$ws.Range("H51").Value = 1091.5
$ws.Range("J51").Value = 1499.6666
$ws.Range("L51").Value = 4498.9998
$ws.Range("N51").Value = -5418.9998
$ws.Range("H55").Value = 1515
$ws.Range("I55").Value = 1072.5
$ws.Range("K55").Value = 3217.5
$ws.Range("M55").Value = -3040.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4749.75
$ws.Range("I70").Value = 4499.5
$ws.Range("K70").Value = 4499.5
$ws.Range("M70").Value = -4229.5
$ws.Range("H73").Value = 4749.75
$ws.Range("I73").Value = 4499.5
$ws.Range("K73").Value = 4499.5
$ws.Range("M73").Value = -3563.5
$ws.Range("H129").Value = 75000
$ws.Range("J129").Value = 75000
$ws.Range("L129").Value = 75000
$ws.Range("N129").Value = -85000
$ws.Range("H132").Value = 3987.2
$ws.Range("I132").Value = 3421.625
$ws.Range("K132").Value = 10264.875
$ws.Range("M132").Value = -7734.875

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 10
$ws.Range("J2").Value = 10
$ws.Range("L2").Value = 10
$ws.Range("N2").Value = -234
$ws.Range("H68").Value = 2221.889
$ws.Range("I68").Value = 2285.4285
$ws.Range("J68").Value = 1999.5
$ws.Range("K68").Value = 2285.4285
$ws.Range("L68").Value = 1999.5
$ws.Range("M68").Value = -1536.4285
$ws.Range("N68").Value = -3497.5
$ws.Range("H71").Value = 2221.889
$ws.Range("I71").Value = 2285.4285
$ws.Range("J71").Value = 1999.5
$ws.Range("K71").Value = 11427.1425
$ws.Range("L71").Value = 9997.5
$ws.Range("M71").Value = -7683.1425
$ws.Range("N71").Value = -17485.5
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H116").Value = 100000
$ws.Range("J116").Value = 100000
$ws.Range("L116").Value = 100000
$ws.Range("N116").Value = -109178
$ws.Range("H120").Value = 200000
$ws.Range("J120").Value = 200000
$ws.Range("L120").Value = 200000
$ws.Range("N120").Value = -209676
$ws.Range("H125").Value = 25000
$ws.Range("J125").Value = 25000
$ws.Range("L125").Value = 25000
$ws.Range("N125").Value = -34840
$ws.Range("H128").Value = 20000
$ws.Range("J128").Value = 20000
$ws.Range("L128").Value = 20000
$ws.Range("N128").Value = -29960
$ws.Range("H132").Value = 5418
$ws.Range("I132").Value = 5501.6
$ws.Range("K132").Value = 16504.8
$ws.Range("M132").Value = -13974.8

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 10508
$ws.Range("I30").Value = 10508
$ws.Range("K30").Value = 10508
$ws.Range("M30").Value = -10401
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H109").Value = 33333
$ws.Range("J109").Value = 33333
$ws.Range("L109").Value = 33333
$ws.Range("N109").Value = -36107
$ws.Range("H110").Value = 14999
$ws.Range("J110").Value = 14999
$ws.Range("L110").Value = 14999
$ws.Range("N110").Value = -23179
$ws.Range("H114").Value = 48333.332
$ws.Range("J114").Value = 45000
$ws.Range("L114").Value = 45000
$ws.Range("N114").Value = -53678
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
